$wb = $excel.ActiveWorkbook

# --- Portfolio sheet: alternative method return for Domestic Equities (B3) ---
$wsPortfolio = $wb.Worksheets.Item("Portfolio")
$wsPortfolio.Activate()
$wsPortfolio.Range("B3").Value = 0.05
$wsPortfolio.Range("B3").Select() | Out-Null

# --- VCV Matrix sheet: no data changes; selection/active state untouched ---

# --- Views sheet becomes the active sheet, with H11 as the selected cell ---
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Activate()
$wsViews.Range("H11").Select() | Out-Null
